$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value2 = 4.2
$ws.Range("L2").Value2 = 2.2
$ws.Range("R2").Value2 = 5.5
$ws.Range("S2").Value2 = 7.6
$ws.Range("V2").Value2 = 17.5
$ws.Range("X2").Value2 = 7.5
$ws.Range("Z2").Value2 = 18.5
$ws.Range("AA2").Value2 = 110
$ws.Range("AB2").Value2 = 9.5
$ws.Range("AC2").Value2 = 22
$ws.Range("AF2").Value2 = 45
$ws.Range("AG2").Value2 = 60
$ws.Range("J3").Value2 = 1.36
$ws.Range("AI3").Value2 = 1.07
$ws.Range("I4").Value2 = 3.1
$ws.Range("J4").Value2 = 1.36
$ws.Range("R4").Value2 = 7.5
$ws.Range("T4").Value2 = 9.5
$ws.Range("X4").Value2 = 9
$ws.Range("AD4").Value2 = 11
$ws.Range("AI4").Value2 = 1.07
$ws.Range("K5").Value2 = 2.5
$ws.Range("L5").Value2 = 2.63
$ws.Range("M5").Value2 = 1.5
$ws.Range("AI5").Value2 = 1.1
$ws.Range("AJ5").Value2 = 7
$ws.Range("J6").Value2 = 1.29
$ws.Range("AI6").Value2 = 1.05
$ws.Range("H7").Value2 = 3.3
$ws.Range("I7").Value2 = 3.2
$ws.Range("J7").Value2 = 1.44
$ws.Range("K7").Value2 = 2.75
$ws.Range("L7").Value2 = 2.38
$ws.Range("M7").Value2 = 1.57
$ws.Range("X7").Value2 = 7.5
$ws.Range("Y7").Value2 = 6
$ws.Range("Z7").Value2 = 17
$ws.Range("AB7").Value2 = 8
$ws.Range("AC7").Value2 = 15
$ws.Range("AD7").Value2 = 12
$ws.Range("AH7").Value2 = 451
$ws.Range("AI7").Value2 = 1.08
$ws.Range("AJ7").Value2 = 7.5
$ws.Range("J8").Value2 = 1.44
$ws.Range("K8").Value2 = 2.75
$ws.Range("AI8").Value2 = 1.08
$ws.Range("G11").Value2 = 4.33
$ws.Range("H11").Value2 = 3.8
$ws.Range("I11").Value2 = 1.73
$ws.Range("J11").Value2 = 1.25
$ws.Range("K11").Value2 = 4
$ws.Range("P11").Value2 = 1.8
$ws.Range("Q11").Value2 = 1.91
$ws.Range("R11").Value2 = 13
$ws.Range("S11").Value2 = 23
$ws.Range("T11").Value2 = 15
$ws.Range("U11").Value2 = 51
$ws.Range("Y11").Value2 = 7.5
$ws.Range("AE11").Value2 = 13
$ws.Range("AF11").Value2 = 13
$ws.Range("AH11").Value2 = 251
$ws.Range("L12").Value2 = 1.83
$ws.Range("M12").Value2 = 2.03
$ws.Range("L15").Value2 = 2.03
$ws.Range("M15").Value2 = 1.83
$ws.Range("H16").Value2 = 4
$ws.Range("N16").Value2 = 1.3
$ws.Range("X16").Value2 = 13
$ws.Range("AB16").Value2 = 8
$ws.Range("AC16").Value2 = 8.5
$ws.Range("AG16").Value2 = 23
$ws.Range("AH16").Value2 = 201
$ws.Range("AI16").Value2 = 1.04
$ws.Range("AJ16").Value2 = 13
$ws.Range("N17").Value2 = 1.41
$ws.Range("O17").Value2 = 2.62
$ws.Range("H18").Value2 = 3.6
$ws.Range("J18").Value2 = 1.33
$ws.Range("K18").Value2 = 3.25
$ws.Range("L18").Value2 = 2.08
$ws.Range("M18").Value2 = 1.73
$ws.Range("N18").Value2 = 1.41
$ws.Range("O18").Value2 = 2.62
$ws.Range("P18").Value2 = 1.91
$ws.Range("Q18").Value2 = 1.8
$ws.Range("R18").Value2 = 6.5
$ws.Range("Z18").Value2 = 17
$ws.Range("AA18").Value2 = 51
$ws.Range("AI18").Value2 = 1.07
$ws.Range("AJ18").Value2 = 9
$ws.Range("G19").Value2 = 1.42
$ws.Range("H19").Value2 = 4.3
$ws.Range("I19").Value2 = 5.7
$ws.Range("L19").Value2 = 1.52
$ws.Range("M19").Value2 = 2.37
$ws.Range("P19").Value2 = 1.69
$ws.Range("Q19").Value2 = 2.04
$ws.Range("R19").Value2 = 7.6
$ws.Range("S19").Value2 = 6.9
$ws.Range("T19").Value2 = 7
$ws.Range("V19").Value2 = 9
$ws.Range("W19").Value2 = 17
$ws.Range("X19").Value2 = 15
$ws.Range("Y19").Value2 = 7.7
$ws.Range("AA19").Value2 = 45
$ws.Range("AB19").Value2 = 16
$ws.Range("AC19").Value2 = 30
$ws.Range("AD19").Value2 = 15
$ws.Range("AE19").Value2 = 80
$ws.Range("AG19").Value2 = 37
$ws.Range("AH19").Value2 = 250
$ws.Range("G20").Value2 = 1.1
$ws.Range("H20").Value2 = 6.7
$ws.Range("I20").Value2 = 18
$ws.Range("L20").Value2 = 1.3
$ws.Range("M20").Value2 = 3.2
$ws.Range("P20").Value2 = 2.12
$ws.Range("Q20").Value2 = 1.65
$ws.Range("S20").Value2 = 6
$ws.Range("U20").Value2 = 5.7
$ws.Range("Z20").Value2 = 27
$ws.Range("AA20").Value2 = 100
$ws.Range("AB20").Value2 = 55
$ws.Range("AC20").Value2 = 175
$ws.Range("AD20").Value2 = 55
$ws.Range("AF20").Value2 = 250
$ws.Range("AG20").Value2 = 120
$ws.Range("G21").Value2 = 5.3
$ws.Range("H21").Value2 = 4
$ws.Range("I21").Value2 = 1.5
$ws.Range("J21").Value2 = 1.16
$ws.Range("K21").Value2 = 4.5
$ws.Range("L21").Value2 = 1.52
$ws.Range("M21").Value2 = 2.22
$ws.Range("P21").Value2 = 1.65
$ws.Range("Q21").Value2 = 2.11
$ws.Range("R21").Value2 = 15.5
$ws.Range("T21").Value2 = 14
$ws.Range("U21").Value2 = 80
$ws.Range("X21").Value2 = 14.5
$ws.Range("Y21").Value2 = 7.2
$ws.Range("Z21").Value2 = 12
$ws.Range("AA21").Value2 = 40
$ws.Range("AC21").Value2 = 7.2
$ws.Range("AD21").Value2 = 6.9
$ws.Range("AE21").Value2 = 9.5
$ws.Range("P27").Value2 = 2.28
$ws.Range("Q27").Value2 = 1.56
$ws.Range("G28").Value2 = 1.5
$ws.Range("H28").Value2 = 4.05
$ws.Range("I28").Value2 = 5.1
$ws.Range("J28").Value2 = 1.18
$ws.Range("K28").Value2 = 4.6
$ws.Range("L28").Value2 = 1.55
$ws.Range("M28").Value2 = 2.15
$ws.Range("P28").Value2 = 1.72
$ws.Range("Q28").Value2 = 2.05
$ws.Range("R28").Value2 = 7.1
$ws.Range("S28").Value2 = 6.9
$ws.Range("T28").Value2 = 6.9
$ws.Range("U28").Value2 = 9.25
$ws.Range("V28").Value2 = 9.5
$ws.Range("X28").Value2 = 13.5
$ws.Range("Z28").Value2 = 12.5
$ws.Range("AA28").Value2 = 45
$ws.Range("AB28").Value2 = 14
$ws.Range("AC28").Value2 = 27
$ws.Range("AD28").Value2 = 14
$ws.Range("AE28").Value2 = 70
$ws.Range("AF28").Value2 = 37
$ws.Range("AG28").Value2 = 35
$ws.Range("AH28").Value2 = 250
$ws.Range("G29").Value2 = 1.28
$ws.Range("H29").Value2 = 4.05
$ws.Range("I29").Value2 = 12.5
$ws.Range("J29").Value2 = 1.26
$ws.Range("K29").Value2 = 3.65
$ws.Range("L29").Value2 = 1.78
$ws.Range("M29").Value2 = 1.83
$ws.Range("N29").Value2 = 1.39
$ws.Range("O29").Value2 = 2.42
$ws.Range("P29").Value2 = 2.39
$ws.Range("Q29").Value2 = 1.54
$ws.Range("R29").Value2 = 4.9
$ws.Range("S29").Value2 = 4.75
$ws.Range("T29").Value2 = 7.1
$ws.Range("U29").Value2 = 6.4
$ws.Range("V29").Value2 = 9.5
$ws.Range("W29").Value2 = 25
$ws.Range("Y29").Value2 = 7.4
$ws.Range("Z29").Value2 = 19.5
$ws.Range("AA29").Value2 = 100
$ws.Range("AB29").Value2 = 22
$ws.Range("AD29").Value2 = 32
$ws.Range("AE29").Value2 = 450
$ws.Range("AF29").Value2 = 175
$ws.Range("AG29").Value2 = 110
$ws.Range("J31").Value2 = 1.29
$ws.Range("K31").Value2 = 3.5
$ws.Range("L31").Value2 = 1.93
$ws.Range("M31").Value2 = 1.93
$ws.Range("G33").Value2 = 2.67
$ws.Range("H33").Value2 = 3
$ws.Range("J33").Value2 = 1.39
$ws.Range("K33").Value2 = 2.77
$ws.Range("R33").Value2 = 7.3
$ws.Range("S33").Value2 = 12.5
$ws.Range("Y33").Value2 = 5.8
$ws.Range("AB33").Value2 = 7.8
$ws.Range("AF33").Value2 = 23
$ws.Range("AH33").Value2 = 600
